$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1053.2142
$ws.Range("I12").Value = 2240.1667
$ws.Range("J12").Value = 163
$ws.Range("K12").Value = 2240.1667
$ws.Range("L12").Value = 163
$ws.Range("M12").Value = -2070.1667
$ws.Range("N12").Value = -503

$ws.Range("H40").Value = 1468.3334
$ws.Range("I40").Value = 780.25
$ws.Range("J40").Value = 1718.5454
$ws.Range("K40").Value = 780.25
$ws.Range("L40").Value = 1718.5454
$ws.Range("M40").Value = -605.25
$ws.Range("N40").Value = -2068.5454

$ws.Range("H76").Value = 4632788
$ws.Range("I76").Value = 3424.875
$ws.Range("J76").Value = 13891514
$ws.Range("K76").Value = 3424.875
$ws.Range("L76").Value = 13891514
$ws.Range("M76").Value = -3109.875
$ws.Range("N76").Value = -13892144

$ws.Range("H79").Value = 4632788
$ws.Range("I79").Value = 3424.875
$ws.Range("J79").Value = 13891514
$ws.Range("K79").Value = 3424.875
$ws.Range("L79").Value = 13891514
$ws.Range("M79").Value = -2332.875
$ws.Range("N79").Value = -13893698

$ws.Range("H113").Value = 200005400
$ws.Range("I113").Value = 333334340
$ws.Range("J113").Value = 12000
$ws.Range("K113").Value = 333334340
$ws.Range("L113").Value = 12000
$ws.Range("M113").Value = -333331086
$ws.Range("N113").Value = -18508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25001488
$ws.Range("I74").Value = 35714820
$ws.Range("J74").Value = 3719.8333
$ws.Range("K74").Value = 35714820
$ws.Range("L74").Value = 3719.8333
$ws.Range("M74").Value = -35713946
$ws.Range("N74").Value = -5467.8333

$ws.Range("H77").Value = 25001488
$ws.Range("I77").Value = 35714820
$ws.Range("J77").Value = 3719.8333
$ws.Range("K77").Value = 178574100
$ws.Range("L77").Value = 18599.1665
$ws.Range("M77").Value = -178569732
$ws.Range("N77").Value = -27335.1665

$ws.Range("H110").Value = 994
$ws.Range("I110").Value = 851.73334
$ws.Range("J110").Value = 1527.5
$ws.Range("K110").Value = 851.73334
$ws.Range("L110").Value = 1527.5
$ws.Range("M110").Value = 1193.26666
$ws.Range("N110").Value = -5617.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 12938
$ws.Range("J75").Value = 9000
$ws.Range("L75").Value = 9000
$ws.Range("N75").Value = -10872

$ws.Range("H78").Value = 12938
$ws.Range("J78").Value = 9000
$ws.Range("L78").Value = 27000
$ws.Range("N78").Value = -36360

$ws.Range("H105").Value = 1853617.4
$ws.Range("I105").Value = 1533.5
$ws.Range("J105").Value = 2633442
$ws.Range("K105").Value = 1533.5
$ws.Range("L105").Value = 2633442
$ws.Range("M105").Value = 213.5
$ws.Range("N105").Value = -2636936

$ws.Range("H134").Value = 3791.6287
$ws.Range("I134").Value = 3791.6287
$ws.Range("K134").Value = 11374.8861
$ws.Range("M134").Value = -8839.8861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2821.4424
$ws.Range("I31").Value = 1413.4054
$ws.Range("K31").Value = 1413.4054
$ws.Range("M31").Value = -1118.4054

$ws.Range("H34").Value = 2821.4424
$ws.Range("I34").Value = 1413.4054
$ws.Range("K34").Value = 1413.4054
$ws.Range("M34").Value = -1211.4054

$ws.Range("H62").Value = 6701.2
$ws.Range("I62").Value = 6166.6665
$ws.Range("K62").Value = 6166.6665
$ws.Range("M62").Value = -5542.6665

$ws.Range("H65").Value = 6701.2
$ws.Range("I65").Value = 6166.6665
$ws.Range("K65").Value = 30833.3325
$ws.Range("M65").Value = -27713.3325

$ws.Range("H107").Value = 1723.0476
$ws.Range("I107").Value = 581.8889
$ws.Range("J107").Value = 2578.9167
$ws.Range("K107").Value = 581.8889
$ws.Range("L107").Value = 2578.9167
$ws.Range("M107").Value = 1338.1111
$ws.Range("N107").Value = -6418.9167

$ws.Range("H122").Value = 1083.9524
$ws.Range("I122").Value = 980.86664
$ws.Range("J122").Value = 1341.6666
$ws.Range("K122").Value = 2942.59992
$ws.Range("L122").Value = 4024.9998
$ws.Range("M122").Value = -492.5999199999997
$ws.Range("N122").Value = -8924.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 834.5
$ws.Range("J24").Value = 834.5
$ws.Range("L24").Value = 2503.5
$ws.Range("N24").Value = -2963.5

$ws.Range("H75").Value = 1626.25
$ws.Range("I75").Value = 1356.5
$ws.Range("J75").Value = 1716.1666
$ws.Range("K75").Value = 4069.5
$ws.Range("L75").Value = 5148.4998
$ws.Range("M75").Value = -3071.5
$ws.Range("N75").Value = -7144.4998

$ws.Range("H78").Value = 1626.25
$ws.Range("I78").Value = 1356.5
$ws.Range("J78").Value = 1716.1666
$ws.Range("K78").Value = 12208.5
$ws.Range("L78").Value = 15445.4994
$ws.Range("M78").Value = -7216.5
$ws.Range("N78").Value = -25429.4994

$ws.Range("H114").Value = 1509.625
$ws.Range("I114").Value = 2394.5
$ws.Range("J114").Value = 624.75
$ws.Range("K114").Value = 7183.5
$ws.Range("L114").Value = 1874.25
$ws.Range("M114").Value = -3929.5
$ws.Range("N114").Value = -8382.25

$ws.Range("H131").Value = 682.33
$ws.Range("J131").Value = 707.7252999999999
$ws.Range("L131").Value = 2123.1759
$ws.Range("N131").Value = -12203.1759

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3298496.5
$ws.Range("I70").Value = 4395.5
$ws.Range("J70").Value = 6958609
$ws.Range("K70").Value = 4395.5
$ws.Range("L70").Value = 6958609
$ws.Range("M70").Value = -4125.5
$ws.Range("N70").Value = -6959149

$ws.Range("H73").Value = 3298496.5
$ws.Range("I73").Value = 4395.5
$ws.Range("J73").Value = 6958609
$ws.Range("K73").Value = 4395.5
$ws.Range("L73").Value = 6958609
$ws.Range("M73").Value = -3459.5
$ws.Range("N73").Value = -6960481

$ws.Range("H80").Value = 3904.818
$ws.Range("I80").Value = 3487.5
$ws.Range("J80").Value = 4143.2856
$ws.Range("K80").Value = 3487.5
$ws.Range("L80").Value = 4143.2856
$ws.Range("M80").Value = -2489.5
$ws.Range("N80").Value = -6139.2856

$ws.Range("H83").Value = 3904.818
$ws.Range("I83").Value = 3487.5
$ws.Range("J83").Value = 4143.2856
$ws.Range("K83").Value = 17437.5
$ws.Range("L83").Value = 20716.428
$ws.Range("M83").Value = -12445.5
$ws.Range("N83").Value = -30700.428

$ws.Range("H113").Value = 8275.532999999999
$ws.Range("I113").Value = 8979.462
$ws.Range("J113").Value = 3700
$ws.Range("K113").Value = 8979.462
$ws.Range("L113").Value = 3700
$ws.Range("M113").Value = -6809.462
$ws.Range("N113").Value = -8040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 18000
$ws.Range("J104").Value = 18000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -24988

$ws.Range("H113").Value = 1072.5883
$ws.Range("I113").Value = 1195.6666
$ws.Range("J113").Value = 149.5
$ws.Range("K113").Value = 3586.9998
$ws.Range("L113").Value = 448.5
$ws.Range("M113").Value = -1416.9998
$ws.Range("N113").Value = -4788.5

$ws.Range("H132").Value = 1551.5625
$ws.Range("I132").Value = 1050.091
$ws.Range("K132").Value = 3150.273
$ws.Range("M132").Value = -620.2729999999997

$ws.Range("H136").Value = 27167832
$ws.Range("I136").Value = 34409940
$ws.Range("K136").Value = 103229820
$ws.Range("M136").Value = -103227270
